$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-18: coin list rotated up by one; row 6 old entry now ends up at row 18
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'8.398"
$ws.Range("E6").Value = "'1.00%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.882"
$ws.Range("E7").Value = "'-1.77%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.945"
$ws.Range("E8").Value = "'-0.75%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9310"
$ws.Range("E9").Value = "'0.83%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1240"
$ws.Range("E10").Value = "'0.01%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1916"
$ws.Range("E11").Value = "'-0.65%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08951"
$ws.Range("E12").Value = "'-2.55%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03289"
$ws.Range("E13").Value = "'-1.67%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09518"
$ws.Range("E14").Value = "'-0.89%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001381"
$ws.Range("E15").Value = "'-0.20%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006100"
$ws.Range("E16").Value = "'4.63%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.373"
$ws.Range("E17").Value = "'-4.21%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.417"
$ws.Range("E18").Value = "'-0.22%"

# Remaining price/volume refresh across the rest of the sheet
$ws.Range("D2").Value = "'313.95"
$ws.Range("E2").Value = "'0.40%"
$ws.Range("D3").Value = "'37.04"
$ws.Range("E3").Value = "'-1.92%"
$ws.Range("D4").Value = "'5.131"
$ws.Range("E4").Value = "'-0.34%"
$ws.Range("D5").Value = "'0.07934"
$ws.Range("E5").Value = "'0.51%"
$ws.Range("D19").Value = "'0.3467"
$ws.Range("E19").Value = "'0.72%"
$ws.Range("D20").Value = "'6.450"
$ws.Range("E20").Value = "'22.37%"
$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'2.07%"
$ws.Range("D22").Value = "'0.2297"
$ws.Range("E22").Value = "'-11.34%"
$ws.Range("E23").Value = "'-1.13%"
$ws.Range("D24").Value = "'0.001192"
$ws.Range("E24").Value = "'-4.59%"
$ws.Range("D25").Value = "'0.004364"
$ws.Range("E25").Value = "'1.05%"
$ws.Range("D26").Value = "'0.0001319"
$ws.Range("E26").Value = "'8.01%"
$ws.Range("D27").Value = "'0.0003956"
$ws.Range("D39").Value = "'0.02267"
$ws.Range("E39").Value = "'-0.65%"
$ws.Range("D40").Value = "'0.05131"
$ws.Range("E40").Value = "'0.58%"
$ws.Range("D41").Value = "'0.007456"
$ws.Range("E41").Value = "'-0.02%"
$ws.Range("D42").Value = "'0.1379"
$ws.Range("E42").Value = "'1.45%"
$ws.Range("D43").Value = "'0.008472"
$ws.Range("E43").Value = "'-3.72%"
$ws.Range("D44").Value = "'0.002053"
$ws.Range("E44").Value = "'6.39%"
$ws.Range("D45").Value = "'0.007797"
$ws.Range("E45").Value = "'-9.58%"
$ws.Range("D46").Value = "'0.00006320"
$ws.Range("E46").Value = "'-6.27%"
$ws.Range("D47").Value = "'0.00000000745"
$ws.Range("E47").Value = "'-0.70%"
$ws.Range("D48").Value = "'0.002842"
$ws.Range("E48").Value = "'-15.40%"
$ws.Range("D49").Value = "'0.001679"
$ws.Range("E49").Value = "'39.88%"
$ws.Range("D50").Value = "'0.00002087"
$ws.Range("E50").Value = "'-0.70%"
$ws.Range("D51").Value = "'0.0001988"
$ws.Range("E51").Value = "'-0.70%"
